$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @(44811, 30, 20000, 20000, 20000, 1333)
    3  = @(44827, 20, 20000, 20000, 20000, 1333)
    4  = @(44830, 25, 12000, 12000, 12000, 800)
    5  = @(44837, 80, 16000, 16000, 16000, 1067)
    6  = @(44824, 20, 20000, 20000, 20000, 1333)
    7  = @(44819, 100, 20000, 20000, 20000, 1333)
    8  = @(44839, 80, 16000, 16000, 16000, 1067)
    9  = @(44508, 40, 10000, 10000, 10000, 667)
    10 = @(44518, 50, 10000, 10000, 10000, 667)
    11 = @(44838, 10, 20000, 20000, 20000, 1333)
    12 = @(44826, 50, 20000, 20000, 20000, 1333)
    13 = @(44525, 40, 8000, 8000, 8000, 533)
    14 = @(44755, 50, 20000, 20000, 20000, 1333)
    15 = @(44756, 80, 20000, 20000, 20000, 1333)
    16 = @(44845, 20, 16000, 16000, 16000, 1067)
    17 = @(44757, 30, 20000, 20000, 20000, 1333)
    19 = @(44776, 80, 20000, 20000, 20000, 1333)
    20 = @(44825, 30, 20000, 20000, 20000, 1333)
    21 = @(44812, 80, 20000, 20000, 20000, 1333)
    22 = @(45134, 5, 20000, 20000, 20000, 1333)
    23 = @(44749, 50, 20000, 20000, 20000, 1333)
    24 = @(44767, 50, 20000, 20000, 20000, 1333)
    25 = @(44771, 40, 20000, 20000, 20000, 1333)
    26 = @(44813, 20, 20000, 20000, 20000, 1333)
    27 = @(44769, 50, 20000, 20000, 20000, 1333)
}

foreach ($row in $values.Keys) {
    $v = $values[$row]
    $ws.Cells.Item($row, 4).Value = $v[0]
    $ws.Cells.Item($row, 10).Value = $v[1]
    $ws.Cells.Item($row, 11).Value = $v[2]
    $ws.Cells.Item($row, 12).Value = $v[3]
    $ws.Cells.Item($row, 13).Value = $v[4]
    $ws.Cells.Item($row, 16).Value = $v[5]
}
